# Update specific values in column D (RandomForest imputed results)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    12 = -6.909599999999998
    27 = -9.048099999999998
    32 = -9.180199999999992
    36 = -8.202099999999998
    38 = -8.690199999999997
    46 = -8.274999999999999
    54 = -8.301800000000002
    55 = -8.427299999999995
    56 = -7.7714
    67 = -6.289199999999999
    69 = -7.120099999999995
    72 = -7.423900000000002
    83 = -8.9496
    86 = -7.605299999999994
    91 = -6.707799999999997
    93 = -6.4719
    99 = -7.784699999999999
}

foreach ($row in $updates.Keys) {
    $ws.Range("D$row").Value = $updates[$row]
}
